# Update NATMI LR-pair output (Cxcl16-Cxcr6) with refreshed TPM-derived statistics.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.221056333333333
$ws.Range("H2").Value = 6.663169
$ws.Range("I2").Value = 0.2217545441472213
$ws.Range("J2").Value = 0.2217545441472213
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.292811
$ws.Range("N2").Value = 0.878433
$ws.Range("O2").Value = 0.1205152987902963
$ws.Range("P2").Value = 0.1205152987902963
$ws.Range("Q2").Value = 0.6503497260196667
$ws.Range("R2").Value = 5.853147534177
$ws.Range("S2").Value = 0.02672481514600833
$ws.Range("T2").Value = 0.02672481514600833
$ws.Range("G3").Value = 2.221056333333333
$ws.Range("H3").Value = 6.663169
$ws.Range("I3").Value = 0.2217545441472213
$ws.Range("J3").Value = 0.2217545441472213
$ws.Range("O3").Value = 0.4409170013616456
$ws.Range("P3").Value = 0.4409170013616455
$ws.Range("Q3").Value = 2.379368046308556
$ws.Range("R3").Value = 21.414312416777
$ws.Range("S3").Value = 0.09777534864371146
$ws.Range("T3").Value = 0.09777534864371146
$ws.Range("G4").Value = 2.221056333333333
$ws.Range("H4").Value = 6.663169
$ws.Range("I4").Value = 0.2217545441472213
$ws.Range("J4").Value = 0.2217545441472213
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.065569666666667
$ws.Range("N4").Value = 3.196709
$ws.Range("O4").Value = 0.4385676998480583
$ws.Range("P4").Value = 0.4385676998480582
$ws.Range("Q4").Value = 2.366690256757889
$ws.Range("R4").Value = 21.300212310821
$ws.Range("S4").Value = 0.09725438035750153
$ws.Range("T4").Value = 0.09725438035750153
$ws.Range("I5").Value = 0.7112336178950279
$ws.Range("J5").Value = 0.711233617895028
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.292811
$ws.Range("N5").Value = 0.878433
$ws.Range("O5").Value = 0.1205152987902963
$ws.Range("P5").Value = 0.1205152987902963
$ws.Range("Q5").Value = 2.085867463563333
$ws.Range("R5").Value = 18.77280717207
$ws.Range("S5").Value = 0.08571453197032272
$ws.Range("T5").Value = 0.08571453197032274
$ws.Range("I6").Value = 0.7112336178950279
$ws.Range("J6").Value = 0.711233617895028
$ws.Range("O6").Value = 0.4409170013616456
$ws.Range("P6").Value = 0.4409170013616455
$ws.Range("S6").Value = 0.3135949940698701
$ws.Range("T6").Value = 0.3135949940698701
$ws.Range("I7").Value = 0.7112336178950279
$ws.Range("J7").Value = 0.711233617895028
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.065569666666667
$ws.Range("N7").Value = 3.196709
$ws.Range("O7").Value = 0.4385676998480583
$ws.Range("P7").Value = 0.4385676998480582
$ws.Range("Q7").Value = 7.590688525567779
$ws.Range("R7").Value = 68.31619673011001
$ws.Range("S7").Value = 0.3119240918548352
$ws.Range("T7").Value = 0.3119240918548352
$ws.Range("G8").Value = 0.6711793333333332
$ws.Range("H8").Value = 2.013538
$ws.Range("I8").Value = 0.06701183795775068
$ws.Range("J8").Value = 0.06701183795775068
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.292811
$ws.Range("N8").Value = 0.878433
$ws.Range("O8").Value = 0.1205152987902963
$ws.Range("P8").Value = 0.1205152987902963
$ws.Range("Q8").Value = 0.1965286917726666
$ws.Range("R8").Value = 1.768758225954
$ws.Range("S8").Value = 0.008075951673965243
$ws.Range("T8").Value = 0.008075951673965241
$ws.Range("G9").Value = 0.6711793333333332
$ws.Range("H9").Value = 2.013538
$ws.Range("I9").Value = 0.06701183795775068
$ws.Range("J9").Value = 0.06701183795775068
$ws.Range("O9").Value = 0.4409170013616456
$ws.Range("P9").Value = 0.4409170013616455
$ws.Range("Q9").Value = 0.7190194301282222
$ws.Range("R9").Value = 6.471174871153999
$ws.Range("S9").Value = 0.02954665864806393
$ws.Range("T9").Value = 0.02954665864806392
$ws.Range("G10").Value = 0.6711793333333332
$ws.Range("H10").Value = 2.013538
$ws.Range("I10").Value = 0.06701183795775068
$ws.Range("J10").Value = 0.06701183795775068
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.065569666666667
$ws.Range("N10").Value = 3.196709
$ws.Range("O10").Value = 0.4385676998480583
$ws.Range("P10").Value = 0.4385676998480582
$ws.Range("Q10").Value = 0.7151883384935556
$ws.Range("R10").Value = 6.436695046441999
$ws.Range("S10").Value = 0.02938922763572152
$ws.Range("T10").Value = 0.02938922763572151
